$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reference Boundaries")

# Row 69: A value changes, plus new B/C/D values
$ws.Range("A69").Value = 8649024
$ws.Range("B69").Value = 40
$ws.Range("C69").Value = 31
$ws.Range("D69").Value = 21

# Row 70
$ws.Range("B70").Value = 26
$ws.Range("C70").Value = 33
$ws.Range("D70").Value = 17

# Row 71
$ws.Range("B71").Value = 29
$ws.Range("C71").Value = 53
$ws.Range("D71").Value = 19

# Row 72
$ws.Range("B72").Value = 26
$ws.Range("C72").Value = 45
$ws.Range("D72").Value = 24

# Row 73
$ws.Range("B73").Value = 23
$ws.Range("C73").Value = 45
$ws.Range("D73").Value = 23

# Row 74
$ws.Range("B74").Value = 34
$ws.Range("C74").Value = 31
$ws.Range("D74").Value = 20

# Row 75
$ws.Range("B75").Value = 32
$ws.Range("C75").Value = 54
$ws.Range("D75").Value = 20

# Row 76
$ws.Range("B76").Value = 31
$ws.Range("C76").Value = 43
$ws.Range("D76").Value = 17

# Row 77
$ws.Range("B77").Value = 31
$ws.Range("C77").Value = 55
$ws.Range("D77").Value = 21

# Row 78 (E78 already has "70th Ground Truth Point", unchanged)
$ws.Range("B78").Value = 37
$ws.Range("C78").Value = 45
$ws.Range("D78").Value = 20

# Row 79
$ws.Range("B79").Value = 22
$ws.Range("C79").Value = 30
$ws.Range("D79").Value = 25

# Row 80
$ws.Range("B80").Value = 32
$ws.Range("C80").Value = 46
$ws.Range("D80").Value = 24

# Row 81
$ws.Range("B81").Value = 30
$ws.Range("C81").Value = 52
$ws.Range("D81").Value = 25

# Row 82
$ws.Range("B82").Value = 27
$ws.Range("C82").Value = 51
$ws.Range("D82").Value = 30

# Row 83 - new note
$ws.Range("E83").Value = "Not the type we are looking for"

# Row 84
$ws.Range("B84").Value = 21
$ws.Range("C84").Value = 47
$ws.Range("D84").Value = 24

# Row 85
$ws.Range("B85").Value = 31
$ws.Range("C85").Value = 48
$ws.Range("D85").Value = 25

# Row 86
$ws.Range("B86").Value = 37
$ws.Range("C86").Value = 58
$ws.Range("D86").Value = 25

# Row 87 - new note
$ws.Range("E87").Value = "Not right"

# Row 88 - move "80th Ground Truth Point" note off this row, add B/C/D
$ws.Range("E88").ClearContents()
$ws.Range("B88").Value = 33
$ws.Range("C88").Value = 61
$ws.Range("D88").Value = 25

# Row 89
$ws.Range("B89").Value = 36
$ws.Range("C89").Value = 51
$ws.Range("D89").Value = 27

# Row 90 - add B/C/D and move "80th Ground Truth Point" note here
$ws.Range("B90").Value = 22
$ws.Range("C90").Value = 63
$ws.Range("D90").Value = 18
$ws.Range("E90").Value = "80th Ground Truth Point"

# Update sheet view / selection to match final state
$excel.ActiveWindow.ScrollRow = 82
$ws.Range("B91").Select()
